$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("身体状态")

# Fill in the new day's record (2018/3/18 -> serial 43176)
$ws.Range("A14").Value = 43176
$ws.Range("B14").Value = 3
$ws.Range("C14").Value = 3
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = "睡得晚，眼睛有点痛"

# Leave the active selection where the user ended up
[void]$ws.Range("E30").Select()
